$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing quantity (2) for the "SG90 9g Micro Servo Motor" row.
$ws.Range("B4").Value = 2
